# Fixed data importer by probability
# - Adds a new "probability" column (G) to the Tasks sheet with per-row values
# - Makes the Tasks sheet the active tab/sheet (was Operators)

$wb = $excel.ActiveWorkbook

$wsTasks = $wb.Worksheets.Item("Tasks")

# New "probability" column header
$wsTasks.Range("G1").Value = "probability"

# Per-row probability values for rows 2..19
$probabilities = @(1, 1, 1, 1, 0, 0, 1, 1, 1, 1, 1, 0.25, 0.15, 0, 0, 0, 0, 0.25)

for ($i = 0; $i -lt $probabilities.Length; $i++) {
    $row = $i + 2
    $wsTasks.Cells.Item($row, 7).Value = $probabilities[$i]
}

# Make "Tasks" the active/selected sheet (updates workbookView.activeTab and
# tabSelected on the sheets), and set its selection to G8
$wsTasks.Select() | Out-Null
$wsTasks.Range("G8").Select() | Out-Null
